$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 stops being the last row, so it picks up the "interior row" look
# (top+bottom thin border, style ids 8/9) that rows 4 and 5 already use.
$ws.Range("A4:E4").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# The new row 7 is also an interior row -> same formatting as row 4/5/6.
$ws.Range("A7:E7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 8 becomes the new last row, so it gets the borderless look that row 6
# used to have (style ids 4/5), taken from row 2 which already has it.
$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 7 - new record (values entered in this order: C, A, D, E, B)
$ws.Range("C7").Value = " The guild\'s crew seems to run\nitself ragged these days…"
$ws.Range("A7").Value = "SCRIPT/G01P03A/um1411.ssb"
$ws.Range("D7").Value = " Гильдейские команды работают\nдо изнеможения..."
$ws.Range("E7").Value = " Ãéìûäåêòëéå ëïíàîäú ñàáïóàýó\näï éèîåíïçåîéÿ…"
$ws.Range("B7").Value = 133

# Row 8 - new record (values entered in this order: A, C, D, E, B)
$ws.Range("A8").Value = "SCRIPT/G01P03A/um1604.ssb"
$ws.Range("C8").Value = " I\'ve heard rumors that the\nTime Gears are being sealed away…"
$ws.Range("D8").Value = " Ходят слухи, что где-то сейчас\nзапечатывают Шестерни Времени..."
$ws.Range("E8").Value = " Öïäÿó òìôöé, œóï ãäå-óï òåêœàò\nèàðåœàóúâàýó Šåòóåñîé Âñåíåîé..."
$ws.Range("B8").Value = 114

# Row heights: new rows use the same 43.2pt height as the others
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 43.2

# Scroll so the new rows are visible, with E8 as the active/selected cell
$excel.Goto($ws.Range("A4"))
$ws.Range("E8").Select()
